$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the low-ppm resistor's Digi-Key SKU in Y5 with the new placeholder value.
$ws.Range("Y5").Value = "ffffffffffffffffffffffffffffff"

# Unhide columns B, E, G (20160929_rocketlogger sheet no longer needs these hidden).
$ws.Columns("B").Hidden = $false
$ws.Columns("E").Hidden = $false
$ws.Columns("G").Hidden = $false

# Reset the view: clear the scrolled topLeftCell and move the selection to Y6.
$null = $ws.Range("Y6").Select()
